$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental"): set Value cell (B7) to the literal text "false".
# A plain $ws.Range("B7").Value = "false" assignment gets auto-coerced by
# Excel into a native Boolean (t="b") cell, which is not what the target
# workbook stores (it keeps a shared-string "false"). Route the literal
# through a formula result + paste-values so it lands as text, matching
# the original author's programmatic edit.
$ws.Range("D1").Formula = "=""false"""
$ws.Range("D1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("D1").ClearContents()

# Row 8 ("Date"): refresh the generated timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Row 17 ("Description"): populate the previously-empty description text.
$ws.Range("B17").Value = "Cardiovascular risk categories based on CRF levels"
